# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values re-pulled from source data
$updates = @{
    3  = -4
    4  = -6
    6  = -3
    8  = -7
    11 = 4
    14 = -6
    18 = -4
    22 = 1
    24 = 1
    28 = -5
    29 = -5
    31 = 4
    32 = -5
    35 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
